# Auto-generated edit script
# Applies scheduled-runner market-price refresh to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 500656.5
$ws.Range("J3").Value = 500656.5
$ws.Range("L3").Value = 500656.5
$ws.Range("N3").Value = -500884.5

$ws.Range("H18").Value = 287.4
$ws.Range("I18").Value = 287.4
$ws.Range("K18").Value = 287.4
$ws.Range("M18").Value = -3.399999999999977

$ws.Range("H40").Value = 1781.375
$ws.Range("I40").Value = 1590
$ws.Range("J40").Value = 2100.3333
$ws.Range("K40").Value = 1590
$ws.Range("L40").Value = 2100.3333
$ws.Range("M40").Value = -1415
$ws.Range("N40").Value = -2450.3333

$ws.Range("H102").Value = 500656.5
$ws.Range("J102").Value = 500656.5
$ws.Range("L102").Value = 500656.5
$ws.Range("N102").Value = -507146.5

$ws.Range("H105").Value = 500335.5
$ws.Range("J105").Value = 500335.5
$ws.Range("L105").Value = 500335.5
$ws.Range("N105").Value = -507323.5

$ws.Range("H116").Value = 8655524
$ws.Range("I116").Value = 17303498
$ws.Range("J116").Value = 7550
$ws.Range("K116").Value = 17303498
$ws.Range("L116").Value = 7550
$ws.Range("M116").Value = -17300056
$ws.Range("N116").Value = -14434

$ws.Range("H133").Value = 27000
$ws.Range("J133").Value = 27000
$ws.Range("L133").Value = 27000
$ws.Range("N133").Value = -37120

$ws.Range("H137").Value = 1301.2122
$ws.Range("I137").Value = 729.41174
$ws.Range("J137").Value = 1908.75
$ws.Range("K137").Value = 2188.23522
$ws.Range("L137").Value = 5726.25
$ws.Range("M137").Value = 361.76478
$ws.Range("N137").Value = -10826.25

$ws.Range("H138").Value = 5594418.5
$ws.Range("I138").Value = 948796.5
$ws.Range("J138").Value = 8067088
$ws.Range("K138").Value = 2846389.5
$ws.Range("L138").Value = 24201264
$ws.Range("M138").Value = -2841249.5
$ws.Range("N138").Value = -24211544

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 970.9697
$ws.Range("I2").Value = 595.34784
$ws.Range("K2").Value = 595.34784
$ws.Range("M2").Value = -482.34784

$ws.Range("H45").Value = 2057.5715
$ws.Range("I45").Value = 1408.6154
$ws.Range("J45").Value = 3112.125
$ws.Range("K45").Value = 1408.6154
$ws.Range("L45").Value = 3112.125
$ws.Range("M45").Value = -1031.6154
$ws.Range("N45").Value = -3866.125

$ws.Range("H74").Value = 1454.7037
$ws.Range("I74").Value = 1531
$ws.Range("K74").Value = 1531
$ws.Range("M74").Value = -657

$ws.Range("H77").Value = 1454.7037
$ws.Range("I77").Value = 1531
$ws.Range("K77").Value = 7655
$ws.Range("M77").Value = -3287

$ws.Range("H97").Value = 18519214
$ws.Range("I97").Value = 25641736
$ws.Range("J97").Value = 662.2
$ws.Range("K97").Value = 25641736
$ws.Range("L97").Value = 662.2
$ws.Range("M97").Value = -25641240
$ws.Range("N97").Value = -1654.2

$ws.Range("H102").Value = 1818.3334
$ws.Range("I102").Value = 1477.5
$ws.Range("K102").Value = 1477.5
$ws.Range("M102").Value = 144.5

$ws.Range("H116").Value = 970.9697
$ws.Range("I116").Value = 595.34784
$ws.Range("K116").Value = 595.34784
$ws.Range("M116").Value = 1698.65216

$ws.Range("H123").Value = 33608
$ws.Range("J123").Value = 33608
$ws.Range("L123").Value = 33608
$ws.Range("N123").Value = -43408

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 970.9697
$ws.Range("I3").Value = 595.34784
$ws.Range("K3").Value = 595.34784
$ws.Range("M3").Value = -481.34784

$ws.Range("H86").Value = 15114.125
$ws.Range("I86").Value = 2135.3333
$ws.Range("K86").Value = 2135.3333
$ws.Range("M86").Value = -1012.3333

$ws.Range("H89").Value = 15114.125
$ws.Range("I89").Value = 2135.3333
$ws.Range("K89").Value = 10676.6665
$ws.Range("M89").Value = -5060.666499999999

$ws.Range("H99").Value = 3312.3076
$ws.Range("I99").Value = 1968.091
$ws.Range("J99").Value = 10705.5
$ws.Range("K99").Value = 1968.091
$ws.Range("L99").Value = 10705.5
$ws.Range("M99").Value = -470.0909999999999
$ws.Range("N99").Value = -13701.5

$ws.Range("H118").Value = 27980
$ws.Range("J118").Value = 27980
$ws.Range("L118").Value = 27980
$ws.Range("N118").Value = -31294

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 500643
$ws.Range("J28").Value = 500643
$ws.Range("L28").Value = 500643
$ws.Range("N28").Value = -501133

$ws.Range("H58").Value = 2051.2593
$ws.Range("I58").Value = 939.5
$ws.Range("J58").Value = 3248.5386
$ws.Range("K58").Value = 939.5
$ws.Range("L58").Value = 3248.5386
$ws.Range("M58").Value = -736.5
$ws.Range("N58").Value = -3654.5386

$ws.Range("H107").Value = 861.6316
$ws.Range("I107").Value = 893.4545000000001
$ws.Range("J107").Value = 817.875
$ws.Range("K107").Value = 893.4545000000001
$ws.Range("L107").Value = 817.875
$ws.Range("M107").Value = 1026.5455
$ws.Range("N107").Value = -4657.875

$ws.Range("H132").Value = 2306.8215
$ws.Range("I132").Value = 1360.0952
$ws.Range("J132").Value = 5147
$ws.Range("K132").Value = 4080.2856
$ws.Range("L132").Value = 15441
$ws.Range("M132").Value = -1550.2856
$ws.Range("N132").Value = -20501

$ws.Range("H134").Value = 2354.0881
$ws.Range("I134").Value = 1176.3334
$ws.Range("J134").Value = 5180.7
$ws.Range("K134").Value = 3529.0002
$ws.Range("L134").Value = 15542.1
$ws.Range("M134").Value = -994.0001999999999
$ws.Range("N134").Value = -20612.1

$ws.Range("H136").Value = 2051.2593
$ws.Range("I136").Value = 939.5
$ws.Range("J136").Value = 3248.5386
$ws.Range("K136").Value = 2818.5
$ws.Range("L136").Value = 9745.6158
$ws.Range("M136").Value = -268.5
$ws.Range("N136").Value = -14845.6158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 975.5204
$ws.Range("I68").Value = 757.69354
$ws.Range("J68").Value = 1350.6666
$ws.Range("K68").Value = 2273.08062
$ws.Range("L68").Value = 4051.9998
$ws.Range("M68").Value = -1462.08062
$ws.Range("N68").Value = -5673.9998

$ws.Range("H71").Value = 975.5204
$ws.Range("I71").Value = 757.69354
$ws.Range("J71").Value = 1350.6666
$ws.Range("K71").Value = 6819.24186
$ws.Range("L71").Value = 12155.9994
$ws.Range("M71").Value = -2763.24186
$ws.Range("N71").Value = -20267.9994

$ws.Range("H102").Value = 4029
$ws.Range("J102").Value = 4029
$ws.Range("L102").Value = 12087
$ws.Range("N102").Value = -16955

$ws.Range("H108").Value = 3600
$ws.Range("J108").Value = 6000
$ws.Range("L108").Value = 18000
$ws.Range("N108").Value = -23760

$ws.Range("H117").Value = 1898
$ws.Range("I117").Value = 399
$ws.Range("J117").Value = 2647.5
$ws.Range("K117").Value = 1197
$ws.Range("L117").Value = 7942.5
$ws.Range("M117").Value = 2245
$ws.Range("N117").Value = -14826.5

$ws.Range("H118").Value = 2270.6365
$ws.Range("I118").Value = 1397.4
$ws.Range("J118").Value = 2998.3333
$ws.Range("K118").Value = 4192.200000000001
$ws.Range("L118").Value = 8994.999899999999
$ws.Range("M118").Value = -2949.200000000001
$ws.Range("N118").Value = -11480.9999

$ws.Range("H139").Value = 33336544
$ws.Range("I139").Value = 41668892
$ws.Range("J139").Value = 7151.3335
$ws.Range("K139").Value = 125006676
$ws.Range("L139").Value = 21454.0005
$ws.Range("M139").Value = -125001536
$ws.Range("N139").Value = -31734.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2683.3333
$ws.Range("J80").Value = 2675
$ws.Range("L80").Value = 2675
$ws.Range("N80").Value = -4671

$ws.Range("H83").Value = 2683.3333
$ws.Range("J83").Value = 2675
$ws.Range("L83").Value = 13375
$ws.Range("N83").Value = -23359

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2723.8096
$ws.Range("I100").Value = 2142.8572
$ws.Range("J100").Value = 3014.2856
$ws.Range("K100").Value = 2142.8572
$ws.Range("L100").Value = 3014.2856
$ws.Range("M100").Value = -1601.8572
$ws.Range("N100").Value = -4096.2856

$ws.Range("H132").Value = 4835.9443
$ws.Range("I132").Value = 3578
$ws.Range("K132").Value = 10734
$ws.Range("M132").Value = -8204

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 33711
$ws.Range("J123").Value = 33711
$ws.Range("L123").Value = 33711
$ws.Range("N123").Value = -43511

$ws.Range("H132").Value = 12822375
$ws.Range("I132").Value = 18519994
$ws.Range("K132").Value = 55559982
$ws.Range("M132").Value = -55557452

$ws.Range("H136").Value = 9260400
$ws.Range("I136").Value = 12821110
$ws.Range("K136").Value = 38463330
$ws.Range("M136").Value = -38460780

$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 80000
$ws.Range("J140").Value = 80000
$ws.Range("L140").Value = 80000
$ws.Range("N140").Value = -90360

$ws.Range("H141").Value = 80000
$ws.Range("J141").Value = 80000
$ws.Range("L141").Value = 80000
$ws.Range("N141").Value = -90360
